$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FTNT")

# Row 6: Change in inventories
$ws.Range("B6").Value = -63600000.0
$ws.Range("C6").Value = -42200000.0
$ws.Range("D6").Value = -54300000.0
$ws.Range("E6").Value = -60900000.0
$ws.Range("F6").Value = -39200000.0
$ws.Range("G6").Value = -48500000.0

# Row 8: Change in payables and accrued liability
$ws.Range("B8").Value = 338300000.0
$ws.Range("C8").Value = 433000000.0
$ws.Range("D8").Value = 387000000.0
$ws.Range("E8").Value = 285900000.0
$ws.Range("F8").Value = 172800000.0
$ws.Range("G8").Value = 69900000.0

# Row 29: Capital Stock Change - B29 was blank inline string, now becomes a number
$ws.Range("B29").Value = -165600000.0
